# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates column G ("K") for rows 3-60 on the active worksheet with
# recalculated strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    3  = 0
    4  = 3
    5  = 2
    6  = 1
    7  = 0
    8  = 1
    9  = 0
    10 = 3
    11 = 1
    12 = 3
    13 = 1
    14 = 1
    15 = 1
    16 = 2
    17 = 1
    18 = 0
    19 = 1
    20 = 2
    21 = 0
    22 = 0
    23 = 3
    24 = 1
    25 = 1
    26 = 2
    27 = 5
    28 = 1
    29 = 2
    30 = 1
    31 = 0
    32 = 0
    33 = 5
    34 = 0
    35 = 0
    36 = 1
    37 = 3
    38 = 5
    39 = 1
    40 = 0
    41 = 2
    42 = 4
    43 = 0
    44 = 2
    45 = 2
    46 = 1
    47 = 2
    48 = 4
    49 = 0
    50 = 2
    51 = 0
    52 = 1
    53 = 1
    54 = 2
    55 = 2
    56 = 1
    57 = 1
    58 = 0
    59 = 1
    60 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
